# Split ISIC 05T06 (Mining and extraction of energy producing products)
# into ISIC 05 (Coal mining) and ISIC 06 (Oil and gas extraction).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "OECD Mapping" sheet: insert a new row before the current
#    "D05T06 / ISIC 05T06" row (row 3) and populate the two rows with
#    the split categories: D05 (Coal mining) then D06 (Oil and gas).
# ---------------------------------------------------------------------
$wsMap = $wb.Worksheets.Item("OECD Mapping")

$wsMap.Rows.Item(3).Insert()

$wsMap.Cells.Item(3, 1).Value() = "D05: Coal mining"
$wsMap.Cells.Item(3, 2).Value() = "ISIC 05"

$wsMap.Cells.Item(4, 1).Value() = "D06: Oil and gas extraction"
$wsMap.Cells.Item(4, 2).Value() = "ISIC 06"

# (All formulas on "Cost Breakdowns" that reference 'OECD Mapping'!A.. /
#  'OECD Mapping'!B.. rows below row 3 shift automatically because of the
#  row insert above.)

# ---------------------------------------------------------------------
# 2) "SoTCCbIC" sheet: insert a new column before the current
#    "ISIC 05T06" column (column D) and relabel the two header cells.
# ---------------------------------------------------------------------
$wsShare = $wb.Worksheets.Item("SoTCCbIC")

$wsShare.Columns.Item(4).Insert()

$wsShare.Cells.Item(1, 3).Value() = "ISIC 05"
$wsShare.Cells.Item(1, 4).Value() = "ISIC 06"

$wsShare.Cells.Item(2, 4).Formula() = "=SUMIF('Cost Breakdowns'!`$E`$3:`$E`$49,SoTCCbIC!D`$1,'Cost Breakdowns'!`$C`$3:`$C`$49)"
